# Insert a new data row at row 8 (pushing the existing rows 8..109 down to 9..110),
# then populate it with the new weekly "Haba" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 45169
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 1200
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 6500
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 260
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
